# This script applies the updated cryptocurrency price/volume data
# described by the commit, row by row, to the active worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "51.233.62"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.55%  "

# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.960.43"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +0.99%  "

# Row 4
$ws.Range("E4").Value = "  +0.00%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "381.41"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +0.68%  "

# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "102.55"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.82%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.541"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.13%  "

# Row 8
$ws.Range("E8").Value = "  -0.01%  "

# Row 9
$ws.Range("E9").Value = "  +0.00%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "36.57"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -1.17%  "

# Row 11
$ws.Range("E11").Value = "  -0.03%  "

# Row 12
$ws.Range("E12").Value = "  -0.25%  "

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.426.41"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +0.93%  "

# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "18.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -2.30%  "

# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.41"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +0.20%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "2.952.40"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +0.67%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.987"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +4.28%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "51.170.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.62%  "

# Row 19
$ws.Range("E19").Value = "  -5.62%  "

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "7.14"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -2.78%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "12.55"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -4.34%  "

# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.0₃0955"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.27%  "

# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "68.52"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.20%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "262.22"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.15%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.89"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +2.66%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "8.38"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +13.12%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.69"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +3.61%  "

# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.170"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +1.48%  "

# Row 29
$ws.Range("B29").Value = "Dai"
$ws.Range("C29").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.04%  "

# Row 30
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.112"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +8.83%  "

# Row 31
$ws.Range("B31").Value = "EthereumClassic"
$ws.Range("C31").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "25.71"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -0.48%  "

# Row 32
$ws.Range("B32").Value = "Cosmos"
$ws.Range("C32").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "9.82"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -0.34%  "

# Row 33
$ws.Range("B33").Value = "VeChain"
$ws.Range("C33").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.0458"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +6.41%  "

# Row 34
$ws.Range("B34").Value = "InjectiveProtocol"
$ws.Range("C34").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "33.96"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -0.44%  "

# Row 35
$ws.Range("B35").Value = "Toncoin"
$ws.Range("C35").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.06"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.39%  "

# Row 36
$ws.Range("B36").Value = "OKB"
$ws.Range("C36").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "50.43"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -3.67%  "

# Row 37
$ws.Range("B37").Value = "FirstDigitalUSD"
$ws.Range("C37").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.00"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.03%  "

# Row 38
$ws.Range("B38").Value = "LidoDAOToken"
$ws.Range("C38").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.99"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -0.93%  "

# Row 39
$ws.Range("B39").Value = "Celestia"
$ws.Range("C39").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "16.81"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -2.45%  "

# Row 40
$ws.Range("B40").Value = "Stellar"
$ws.Range("C40").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.116"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +0.54%  "

# Row 41
$ws.Range("E41").Value = "  -2.29%  "

# Row 42
$ws.Range("B42").Value = "ARBITRUM"
$ws.Range("C42").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -2.25%  "

# Row 43
$ws.Range("B43").Value = "Monero"
$ws.Range("C43").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "121.63"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.94%  "

# Row 44
$ws.Range("B44").Value = "EnergySwap"
$ws.Range("C44").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "21.33"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -2.90%  "

# Row 45
$ws.Range("B45").Value = "WEMIXToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "2.05"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -0.98%  "

# Row 46
$ws.Range("B46").Value = "TheGraph"
$ws.Range("C46").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.273"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.75%  "

# Row 47
$ws.Range("B47").Value = "ApeXProtocol"
$ws.Range("C47").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.37"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.56%  "

# Row 48
$ws.Range("B48").Value = "NEARProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "3.24"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +1.63%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.012.47"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -0.34%  "

# Row 50
$ws.Range("B50").Value = "BEAM"
$ws.Range("C50").Value = "https://coinranking.com/coin/cYYMfXF4u+beam-beam"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0345"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +6.83%  "

# Row 51
$ws.Range("B51").Value = "StarkNetToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/e8MSzQnGc+starknettoken-strk"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.17"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +17.70%  "

